$wb = $excel.ActiveWorkbook

# Add a new worksheet after the existing Sheet1/Sheet2 tabs. This sheet is
# used to exercise column addresses beyond 'Z' (i.e. two-letter columns
# like "AA"), matching the "Fix spreadsheet column address parsing" fix.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Big col nums"

# Populate the big-column-number regression cells.
$newSheet.Range("Z1").Value = "Cell Z1"
$newSheet.Range("AA2").Value = "Cell AA2"
